# Update the "Metadata" sheet of the ValueSet workbook:
#   - B5 (Title)  : "Vaccine Contraindication" -> "NG-Imm Vaccine Contraindication VS"
#   - B8 (Date)   : "2025-06-23T13:45:54+01:00" -> "2025-06-24T09:13:37+01:00"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B5").Value = "NG-Imm Vaccine Contraindication VS"
$ws.Range("B8").Value = "2025-06-24T09:13:37+01:00"
